$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 12 that reuses row 11's current formatting (border/font/
#     wrap/alignment, row height), mirroring a "copy row 11 -> paste into row 12"
#     edit in the UI. We copy formats only (values/text are set explicitly below).
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("12").RowHeight = 47.25

# --- Center the "Impacto" column value on row 11
$ws.Range("C11").HorizontalAlignment = -4108

# --- Populate the new row 12 with the new "Bloqueo de cuenta" finding, and fix
#     the ReferenceID text on row 11 (167 -> 144 was wrong, should be 141), in
#     the same order the cells were authored.
$ws.Range("B12").Value = "Bloqueo de cuenta innactivo"
$ws.Range("C12").Value = "Medio"
$ws.Range("D12").Value = "Contador de numero de intentos permitidos"
$ws.Range("E11").Value = "El ReferenceID cambia a de 167 a 141"
$ws.Range("E12").Value = "El usuario deberá vizualizar la cantidad intentos permitidos asi mismo redirigirlo a una pagina donde le notifique que la cuenta fue bloqueda."
$ws.Range("F12").Value = "Error de codigo fuente"
$ws.Range("G12").Value = "Abierta"

# --- Match the view state recorded after the edit
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B13").Select()
